$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 48, shifting the old row 48 (and below) down to 49.
$ws.Rows.Item(48).Insert()

# New row 48 gets a copy of what row 47 currently holds (the "old" weekly entry),
# since the new week's data is written into row 47 in place.
$ws.Range("A48").Value = 2
$ws.Range("B48").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C48").Value = "Coquimbo"
$ws.Range("D48").Value = 44692
$ws.Range("D48").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E48").Value = 4
$ws.Range("F48").Value = 100112026
$ws.Range("G48").Value = "Haba"
$ws.Range("H48").Value = "Sin especificar"
$ws.Range("I48").Value = "Primera"
$ws.Range("J48").Value = 500
$ws.Range("K48").Value = 15000
$ws.Range("L48").Value = 16000
$ws.Range("M48").Value = 15500
$ws.Range("N48").Value = "$/saco 25 kilos"
$ws.Range("O48").Value = "Provincia de Limarí"
$ws.Range("P48").Value = 620
$ws.Range("Q48").Value = 25
$ws.Range("R48").Value = "Hortaliza"

# Update row 47 in place with the new week's data.
$ws.Range("D47").Value = 44706
$ws.Range("K47").Value = 13000
$ws.Range("L47").Value = 14000
$ws.Range("M47").Value = 13500
$ws.Range("P47").Value = 540
